$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (shifts old N:P to O:Q),
# copying the column width from the preceding column (M) as Excel does
# by default when inserting a column.
$mWidth = $ws.Columns("M").ColumnWidth
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $mWidth

# Make "Repayment schedule" the active sheet/tab and set its selection.
$ws.Activate()
$ws.Range("Q7").Select() | Out-Null
